# M7 Deployment deck update
#  - Add "Three Deployment Methods" slide (FTP / Visual Studio Publish / Publish
#    from source control) before the existing "Summary" slide.
#  - Add "Resources" slide (aka.ms/mdha, aka.ms/kudu, aka.ms/azurewebsites)
#    before the existing "End" slide.
#  - Tidy two runs that were split across edits into single runs.

$p = $ppt.ActivePresentation

# "Title and Content" custom layout on the main (first) slide master - the
# same layout already used by the neighbouring "Summary" slide.
$layout = $p.SlideMaster.CustomLayouts.Item(4)

# --- New slide 18: Three Deployment Methods ------------------------------
$deploy = $p.Slides.AddSlide(18, $layout)

$deployTitle = $deploy.Shapes.Item(1)
$deployTitle.TextFrame.TextRange.Text = "Three Deployment Methods"

$deployBody = $deploy.Shapes.Item(2)
$deployBody.TextFrame.TextRange.Text = "FTP`rVisual Studio Publish`rPublish from source control"

# --- New slide 19: Resources ----------------------------------------------
$resources = $p.Slides.AddSlide(19, $layout)

$resourcesTitle = $resources.Shapes.Item(1)
$resourcesTitle.TextFrame.TextRange.Text = "Resources"

$resourcesBody = $resources.Shapes.Item(2)
$resourcesBody.TextFrame.TextRange.Text = "aka.ms/mdha`raka.ms/kudu`raka.ms/azurewebsites"

# --- Slide 4 ("Multi-Device Hybrid Apps Tooling" module intro): merge the
# two runs that make up the title line into a single run. ------------------
$introSlide = $p.Slides.Item(4)
$introShape = $introSlide.Shapes.Item(1)
$introRange = $introShape.TextFrame.TextRange
$introFirstLine = $introRange.Characters(1, 33)
$introFirstLine.Text = "Multi-Device Hybrid Apps Tooling"

# --- Slide 16 ("Email the product team" resources slide): merge the two
# runs "Email " and "the product team" into a single run. -------------------
$emailSlide = $p.Slides.Item(16)
$emailShape = $emailSlide.Shapes.Item(2)
$emailRange = $emailShape.TextFrame.TextRange
$emailPhrase = $emailRange.Characters(87, 22)
$emailPhrase.Text = "Email the product team"
